$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 310
$ws.Range("F7").Value = 140
$ws.Range("F10").Value = 55
$ws.Range("F13").Value = 2863
$ws.Range("F15").Value = 37
$ws.Range("F17").Value = 27
$ws.Range("F20").Value = 21
$ws.Range("F21").Value = 636
$ws.Range("F23").Value = 104
$ws.Range("F25").Value = 32
$ws.Range("F27").Value = 2286
$ws.Range("F28").Value = 4836
$ws.Range("F32").Value = 1251
$ws.Range("F33").Value = 254
$ws.Range("F34").Value = 2167
$ws.Range("F36").Value = 477
$ws.Range("F38").Value = 63
$ws.Range("F39").Value = 143
$ws.Range("F41").Value = 448
$ws.Range("F42").Value = 755
$ws.Range("F46").Value = 444

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 310
$ws.Range("F7").Value = 140
$ws.Range("F10").Value = 55
$ws.Range("F13").Value = 2863
$ws.Range("F15").Value = 37
$ws.Range("F18").Value = 27
$ws.Range("F21").Value = 21
$ws.Range("F22").Value = 636
$ws.Range("F24").Value = 104
$ws.Range("F26").Value = 32
$ws.Range("F28").Value = 2286
$ws.Range("F29").Value = 4836
$ws.Range("F33").Value = 1251
$ws.Range("F34").Value = 254
$ws.Range("F35").Value = 2167
$ws.Range("F37").Value = 477
$ws.Range("F39").Value = 63
$ws.Range("F40").Value = 143
$ws.Range("F42").Value = 448
$ws.Range("F43").Value = 755
$ws.Range("F47").Value = 444
